# Update row 17 (15th data record) of the Binance ETHUSDT kline data with
# the corrected close / volume / quote-asset-volume / trade-count / taker
# figures, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking values are stored as text (to preserve trailing zeros),
# so prefix with a leading apostrophe to force text entry and keep full
# precision instead of Excel auto-converting them to floating point numbers.
$ws.Range("F17").Value = "'142.22000000"
$ws.Range("G17").Value = "'97817.78584000"
$ws.Range("I17").Value = "'13822150.19547320"
$ws.Range("J17").Value = 48956
$ws.Range("K17").Value = "'45497.97767000"
$ws.Range("L17").Value = "'6431745.54879430"
